# Add a "color style" column (column E) to every sheet, one hex color per
# topic, and leave the workbook positioned on the Geografia tab.
#
# Order of operations mirrors how a person would have clicked through the
# workbook tab by tab (Historia -> Entretenimiento -> Deporte -> Geografia),
# which is also the order the new shared strings end up in.

$wb = $excel.ActiveWorkbook

# --- Historia: #daa520 -------------------------------------------------
$wsHistoria = $wb.Worksheets.Item("Historia")
[void]$wsHistoria.Activate()
$wsHistoria.Range("E1").Value = "#daa520"
$wsHistoria.Range("E2").Value = "#daa520"
$wsHistoria.Range("E3").Value = "#daa520"
$wsHistoria.PageSetup.Orientation = 1
[void]$wsHistoria.Range("E3").Select()

# --- Entretenimiento: #db7093 ------------------------------------------
$wsEntretenimiento = $wb.Worksheets.Item("Entretenimiento")
[void]$wsEntretenimiento.Activate()
$wsEntretenimiento.Range("E1").Value = "#db7093"
$wsEntretenimiento.Range("E2").Value = "#db7093"
[void]$wsEntretenimiento.Range("E2").Select()

# --- Deporte: #20b2aa ----------------------------------------------------
$wsDeporte = $wb.Worksheets.Item("Deporte")
[void]$wsDeporte.Activate()
$wsDeporte.Range("E1").Value = "#20b2aa"
$wsDeporte.Range("E2").Value = "#20b2aa"
$wsDeporte.Range("E3").Value = "#20b2aa"
[void]$wsDeporte.Range("E3").Select()

# --- Geografia: #3cb371 (ends as the active sheet) ----------------------
$wsGeografia = $wb.Worksheets.Item("Geografia")
[void]$wsGeografia.Activate()
$wsGeografia.Range("E1").Value = "#3cb371"
$wsGeografia.Range("E2").Value = "#3cb371"
$wsGeografia.Range("E3").Value = "#3cb371"
[void]$wsGeografia.Range("G3").Select()
